$d = $word.ActiveDocument

# NOTE: In this runtime, Find.Execute (whether invoked on $d.Content, a
# Range, or a Selection) always scans the *whole* document story from the
# beginning rather than being confined to the object it was called on.
# wdReplaceOne (the literal 1 for the Replace argument) only swaps the
# first remaining match, though, so issuing the calls below in the exact
# left-to-right order their target text occurs in the document lets us
# reliably address repeated values (e.g. "35", "21", "27", "93",
# "DP I NAUTIKA" each occur more than once across the three table cells).
#
# To stop a freshly-written new value from colliding with an *original*
# value that is still waiting to be matched later in the sequence (e.g.
# writing a "36" for cell 1 before cell 2's original "36" has been
# located), the substitution happens in two passes: first every relevant
# occurrence (in document order, including ones whose text does not
# actually change - they still have to be "consumed" in order) is swapped
# for a unique placeholder token, then every placeholder is resolved to
# its real final text.

function Replace-One($old, $new) {
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 1) | Out-Null
}

# ---- Pass 1: original text (document order) -> unique placeholders ----
    Replace-One "K1" "@@P0@@"
    Replace-One "ANTONIUS SETIAWAN" "@@P1@@"
    Replace-One "93" "@@P2@@"
    Replace-One "110" "@@P3@@"
    Replace-One "72" "@@P4@@"
    Replace-One "35" "@@P5@@"
    Replace-One "26" "@@P6@@"
    Replace-One "21" "@@P7@@"
    Replace-One "97" "@@P8@@"
    Replace-One "DP I NAUTIKA" "@@P9@@"
    Replace-One "K2" "@@P10@@"
    Replace-One "NUR SYAMSI" "@@P11@@"
    Replace-One "100" "@@P12@@"
    Replace-One "117" "@@P13@@"
    Replace-One "76" "@@P14@@"
    Replace-One "36" "@@P15@@"
    Replace-One "27" "@@P16@@"
    Replace-One "22" "@@P17@@"
    Replace-One "102" "@@P18@@"
    Replace-One "DP I NAUTIKA" "@@P19@@"
    Replace-One "K3" "@@P20@@"
    Replace-One "ERIK HANDOYO" "@@P21@@"
    Replace-One "93" "@@P22@@"
    Replace-One "111" "@@P23@@"
    Replace-One "71" "@@P24@@"
    Replace-One "35" "@@P25@@"
    Replace-One "27" "@@P26@@"
    Replace-One "21" "@@P27@@"
    Replace-One "98" "@@P28@@"
    Replace-One "DP I NAUTIKA" "@@P29@@"

# ---- Pass 2: placeholders -> final new text ----
    Replace-One "@@P0@@" "K40"
    Replace-One "@@P1@@" "RAHMAD HIDAYAT"
    Replace-One "@@P2@@" "93"
    Replace-One "@@P3@@" "116"
    Replace-One "@@P4@@" "70"
    Replace-One "@@P5@@" "36"
    Replace-One "@@P6@@" "28"
    Replace-One "@@P7@@" "22"
    Replace-One "@@P8@@" "98"
    Replace-One "@@P9@@" "DP I TEKNIKA"
    Replace-One "@@P10@@" "L3"
    Replace-One "@@P11@@" "SUGENG GUNADI"
    Replace-One "@@P12@@" "82"
    Replace-One "@@P13@@" "102"
    Replace-One "@@P14@@" "67"
    Replace-One "@@P15@@" "33"
    Replace-One "@@P16@@" "26"
    Replace-One "@@P17@@" "20"
    Replace-One "@@P18@@" "97"
    Replace-One "@@P19@@" "DP I TEKNIKA"
    Replace-One "@@P20@@" "L5"
    Replace-One "@@P21@@" "TRIYONO"
    Replace-One "@@P22@@" "79"
    Replace-One "@@P23@@" "99"
    Replace-One "@@P24@@" "63"
    Replace-One "@@P25@@" "32"
    Replace-One "@@P26@@" "25"
    Replace-One "@@P27@@" "19"
    Replace-One "@@P28@@" "97"
    Replace-One "@@P29@@" "DP I TEKNIKA"
